$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45203) for every data
# row (2 through 472). The commit bumps that serial by one day (45203 ->
# 45204) across the whole column range.
$ws.Range("C2:C472").Value = 45204
